$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.Formula = '="60.599.51"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(2, 5).Value = '  +3.09%  '
$cell = $ws.Cells.Item(3, 4)
$cell.Formula = '="2.700.46"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(3, 5).Value = '  +2.68%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$cell = $ws.Cells.Item(5, 4)
$cell.Formula = '="527.19"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(5, 5).Value = '  +1.43%  '
$cell = $ws.Cells.Item(6, 4)
$cell.Formula = '="149.57"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(6, 5).Value = '  +2.48%  '
$cell = $ws.Cells.Item(7, 4)
$cell.Formula = '="0.997"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(7, 5).Value = '  -0.10%  '
$cell = $ws.Cells.Item(8, 4)
$cell.Formula = '="0.578"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(8, 5).Value = '  +1.27%  '
$cell = $ws.Cells.Item(9, 4)
$cell.Formula = '="2.715.47"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(9, 5).Value = '  +2.89%  '
$cell = $ws.Cells.Item(10, 4)
$cell.Formula = '="7.04"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(10, 5).Value = '  +11.18%  '
$cell = $ws.Cells.Item(11, 4)
$cell.Formula = '="0.105"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(11, 5).Value = '  +0.75%  '
$cell = $ws.Cells.Item(12, 4)
$cell.Formula = '="0.341"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(12, 5).Value = '  +1.79%  '
$cell = $ws.Cells.Item(13, 4)
$cell.Formula = '="0.130"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(13, 5).Value = '  +2.44%  '
$cell = $ws.Cells.Item(14, 4)
$cell.Formula = '="3.177.04"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(14, 5).Value = '  +2.65%  '
$cell = $ws.Cells.Item(15, 4)
$cell.Formula = '="60.585.76"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(15, 5).Value = '  +3.01%  '
$cell = $ws.Cells.Item(16, 4)
$cell.Formula = '="21.51"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(16, 5).Value = '  +3.24%  '
$ws.Cells.Item(17, 2).Value = 'ShibaInu'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Cells.Item(17, 4)
$cell.Formula = '="0.0000139"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(17, 5).Value = '  +1.27%  '
$ws.Cells.Item(18, 2).Value = 'WrappedEther'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$cell = $ws.Cells.Item(18, 4)
$cell.Formula = '="2.707.87"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(18, 5).Value = '  +2.54%  '
$cell = $ws.Cells.Item(19, 4)
$cell.Formula = '="346.02"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(19, 5).Value = '  -0.62%  '
$cell = $ws.Cells.Item(20, 4)
$cell.Formula = '="4.51"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(20, 5).Value = '  +1.02%  '
$cell = $ws.Cells.Item(21, 4)
$cell.Formula = '="10.54"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(21, 5).Value = '  +2.67%  '
$cell = $ws.Cells.Item(22, 4)
$cell.Formula = '="6.40"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(22, 5).Value = '  +3.92%  '
$cell = $ws.Cells.Item(23, 4)
$cell.Formula = '="0.998"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(23, 5).Value = '  +0.05%  '
$cell = $ws.Cells.Item(24, 4)
$cell.Formula = '="63.63"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(24, 5).Value = '  +2.91%  '
$cell = $ws.Cells.Item(25, 4)
$cell.Formula = '="0.170"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(25, 5).Value = '  +4.34%  '
$cell = $ws.Cells.Item(26, 4)
$cell.Formula = '="0.419"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(26, 5).Value = '  +1.02%  '
$ws.Cells.Item(27, 5).Value = '  -0.12%  '
$cell = $ws.Cells.Item(28, 4)
$cell.Formula = '="0.0₃0822"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(28, 5).Value = '  +2.53%  '
$cell = $ws.Cells.Item(29, 4)
$cell.Formula = '="7.29"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(29, 5).Value = '  +3.29%  '
$ws.Cells.Item(30, 5).Value = '  +8.30%  '
$cell = $ws.Cells.Item(31, 4)
$cell.Formula = '="0.998"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(31, 5).Value = '  +0.00%  '
$ws.Cells.Item(32, 5).Value = '  +1.72%  '
$cell = $ws.Cells.Item(33, 4)
$cell.Formula = '="19.08"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(33, 5).Value = '  +1.15%  '
$cell = $ws.Cells.Item(34, 4)
$cell.Formula = '="150.74"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(34, 5).Value = '  +1.00%  '
$cell = $ws.Cells.Item(35, 4)
$cell.Formula = '="4.27"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(35, 5).Value = '  +6.42%  '
$cell = $ws.Cells.Item(36, 4)
$cell.Formula = '="1.22"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(36, 5).Value = '  +6.35%  '
$cell = $ws.Cells.Item(37, 4)
$cell.Formula = '="0.922"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(37, 5).Value = '  -5.02%  '
$cell = $ws.Cells.Item(38, 4)
$cell.Formula = '="0.903"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(38, 5).Value = '  +7.01%  '
$cell = $ws.Cells.Item(39, 4)
$cell.Formula = '="1.53"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(39, 5).Value = '  +7.70%  '
$cell = $ws.Cells.Item(40, 4)
$cell.Formula = '="37.35"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(40, 5).Value = '  +2.18%  '
$cell = $ws.Cells.Item(41, 4)
$cell.Formula = '="3.65"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(41, 5).Value = '  +0.55%  '
$cell = $ws.Cells.Item(42, 4)
$cell.Formula = '="0.635"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(42, 5).Value = '  +5.73%  '
$cell = $ws.Cells.Item(43, 4)
$cell.Formula = '="279.66"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(43, 5).Value = '  +0.25%  '
$cell = $ws.Cells.Item(44, 4)
$cell.Formula = '="20.11"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(44, 5).Value = '  +2.31%  '
$cell = $ws.Cells.Item(45, 4)
$cell.Formula = '="0.996"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(45, 5).Value = '  -0.17%  '
$cell = $ws.Cells.Item(46, 4)
$cell.Formula = '="0.0986"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(46, 5).Value = '  +0.21%  '
$ws.Cells.Item(47, 5).Value = '  +6.78%  '
$cell = $ws.Cells.Item(48, 4)
$cell.Formula = '="0.0545"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(48, 5).Value = '  +3.76%  '
$cell = $ws.Cells.Item(49, 4)
$cell.Formula = '="2.098.49"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(49, 5).Value = '  +1.17%  '
$ws.Cells.Item(50, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$cell = $ws.Cells.Item(50, 4)
$cell.Formula = '="10.55"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(50, 5).Value = '  +2.44%  '
$ws.Cells.Item(51, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Cells.Item(51, 4)
$cell.Formula = '="19.48"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(51, 5).Value = '  +4.57%  '

$excel.CutCopyMode = 0
